$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("files")

# Row 3 (cm.xpt): add Status/Who -> InProgress / JJ
$ws.Range("C3").Value = "InProgress"
$ws.Range("D3").Value = "JJ"

# Row 15 (suppds.xpt): status TobeDone -> InProgress (Who stays Cindy)
$ws.Range("C15").Value = "InProgress"

# Row 17 (sv.xpt): add Status/Who -> InProgress / Jessica
$ws.Range("C17").Value = "InProgress"
$ws.Range("D17").Value = "Jessica"

# Row 23 (vs.xpt): status InProgress -> done (Who stays Jessica)
$ws.Range("C23").Value = "done"

# Update selection/view to match final state
$ws.Range("C18").Select()
